$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.950.62"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.818.55"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D5").Value = "309.94"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.4657"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "0.3664"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("D9").Value = "0.07349"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "0.8728"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "20.26"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "1.846.76"
$ws.Range("E12").Value = "  +6.72%  "
$ws.Range("D13").Value = "5.412"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "0.07106"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "6.515"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "91.45"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'0.000008715"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "14.66"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "26.971.54"
$ws.Range("D22").Value = "5.294"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D24").Value = "2.050.88"
$ws.Range("E24").Value = "  +4.45%  "
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "150.92"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").Value = "18.37"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'2.150"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "5.256"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "117.35"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").Value = "0.08894"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "0.7587"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "1.161"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "'4.500"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").Value = "'0.05300"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "0.01946"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "2.969"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").Value = "7.189"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "0.5294"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "2.348"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "8.436"
$ws.Range("D46").Value = "0.4876"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").Value = "10.49"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "1.665"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "103.43"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "0.06297"
$ws.Range("E51").Value = "  +0.11%  "
